$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks (cols A:B, E, J:Z get very slightly wider/narrower) ---
$ws.Columns.Item(1).Resize(1,2).ColumnWidth = 9.67   # A:B -> stored width 10.5
$ws.Range("E:E").ColumnWidth = 14.67                  # E -> stored width 15.5
$ws.Range("J1:Z1").EntireColumn.ColumnWidth = 7.67    # J:Z -> stored width 8.5

# --- New column J ("empty" column header) ---
$ws.Range("B3").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("J3").Value = "empty"
$ws.Range("J3").HorizontalAlignment = 1

# --- A4 gets the same "0.00" number style A3/A9 used to have ---
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)

# --- A3 and A9's old formatted (but empty) cells go away entirely ---
$ws.Range("A3").Clear()
$ws.Range("A9").Clear()

# --- A5 gets a brand-new currency-like number format ---
$ws.Range("A5").NumberFormat = "#,##0.00 ""€"""

# --- Selection moves from A9 to A5 ---
$ws.Range("A5").Select()
